$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -20.06446253974092
$ws.Range("C2").Value = 2.589187502720277
$ws.Range("D2").Value = -20.06446253974092
$ws.Range("E2").Value = -20.06446253974092
$ws.Range("F2").Value = -20.06446253974092
$ws.Range("G2").Value = -20.06446253974092
$ws.Range("H2").Value = -20.06446253974092
$ws.Range("I2").Value = -20.06446253974092
$ws.Range("J2").Value = -20.06446253974092
$ws.Range("K2").Value = -20.06446253974092

$ws.Range("B3").Value = -20.06446253974092
$ws.Range("C3").Value = -20.06446253974092
$ws.Range("D3").Value = -20.06446253974092
$ws.Range("E3").Value = -20.06446253974092
$ws.Range("F3").Value = -20.06446253974092
$ws.Range("G3").Value = -20.06446253974092
$ws.Range("H3").Value = -20.06446253974092
$ws.Range("I3").Value = 2.378933559256486
$ws.Range("J3").Value = -20.06446253974092
$ws.Range("K3").Value = -20.06446253974092

$ws.Range("B4").Value = -20.06446253974092
$ws.Range("C4").Value = 1.978892221641074
$ws.Range("D4").Value = 2.974770634834108
$ws.Range("E4").Value = -20.06446253974092
$ws.Range("F4").Value = 2.47905763226863
$ws.Range("G4").Value = -20.06446253974092
$ws.Range("H4").Value = 1.830432547207042
$ws.Range("I4").Value = -20.06446253974092
$ws.Range("J4").Value = 2.374966873696874
$ws.Range("K4").Value = -20.06446253974092

$ws.Range("B5").Value = -20.06446253974092
$ws.Range("C5").Value = 0.8491639327310042
$ws.Range("D5").Value = -20.06446253974092
$ws.Range("E5").Value = -20.06446253974092
$ws.Range("F5").Value = -20.06446253974092
$ws.Range("G5").Value = 2.132549226907248
$ws.Range("H5").Value = -20.06446253974092
$ws.Range("I5").Value = -20.06446253974092
$ws.Range("J5").Value = -20.06446253974092
$ws.Range("K5").Value = -20.06446253974092

$ws.Range("B6").Value = -20.06446253974092
$ws.Range("C6").Value = -20.06446253974092
$ws.Range("D6").Value = -20.06446253974092
$ws.Range("E6").Value = -20.06446253974092
$ws.Range("F6").Value = -20.06446253974092
$ws.Range("G6").Value = -20.06446253974092
$ws.Range("H6").Value = -20.06446253974092
$ws.Range("I6").Value = -20.06446253974092
$ws.Range("J6").Value = -20.06446253974092
$ws.Range("K6").Value = -20.06446253974092

$ws.Range("B7").Value = 4.321926844935577
$ws.Range("C7").Value = -20.06446253974092
$ws.Range("D7").Value = -20.06446253974092
$ws.Range("E7").Value = -20.06446253974092
$ws.Range("F7").Value = -20.06446253974092
$ws.Range("G7").Value = -20.06446253974092
$ws.Range("H7").Value = -20.06446253974092
$ws.Range("I7").Value = -20.06446253974092
$ws.Range("J7").Value = -20.06446253974092
$ws.Range("K7").Value = -20.06446253974092

$ws.Range("B8").Value = -20.06446253974092
$ws.Range("C8").Value = -20.06446253974092
$ws.Range("D8").Value = -20.06446253974092
$ws.Range("E8").Value = 2.907753198756587
$ws.Range("F8").Value = -20.06446253974092
$ws.Range("G8").Value = -20.06446253974092
$ws.Range("H8").Value = -20.06446253974092
$ws.Range("I8").Value = -20.06446253974092
$ws.Range("J8").Value = -20.06446253974092
$ws.Range("K8").Value = -20.06446253974092

$ws.Range("B9").Value = -20.06446253974092
$ws.Range("C9").Value = -20.06446253974092
$ws.Range("D9").Value = -20.06446253974092
$ws.Range("E9").Value = -20.06446253974092
$ws.Range("F9").Value = -20.06446253974092
$ws.Range("G9").Value = -20.06446253974092
$ws.Range("H9").Value = -20.06446253974092
$ws.Range("I9").Value = -20.06446253974092
$ws.Range("J9").Value = -20.06446253974092
$ws.Range("K9").Value = -20.06446253974092

$ws.Range("B10").Value = -20.06446253974092
$ws.Range("C10").Value = -20.06446253974092
$ws.Range("D10").Value = -20.06446253974092
$ws.Range("E10").Value = -20.06446253974092
$ws.Range("F10").Value = -20.06446253974092
$ws.Range("G10").Value = -20.06446253974092
$ws.Range("H10").Value = -20.06446253974092
$ws.Range("I10").Value = 1.545286583574901
$ws.Range("J10").Value = -20.06446253974092
$ws.Range("K10").Value = 2.207736942485943

$ws.Range("B11").Value = -20.06446253974092
$ws.Range("C11").Value = -20.06446253974092
$ws.Range("D11").Value = -20.06446253974092
$ws.Range("E11").Value = 2.006503813758147
$ws.Range("F11").Value = -20.06446253974092
$ws.Range("G11").Value = 2.584702619030688
$ws.Range("H11").Value = -20.06446253974092
$ws.Range("I11").Value = -20.06446253974092
$ws.Range("J11").Value = -20.06446253974092
$ws.Range("K11").Value = 1.374307489213827

$ws.Range("B12").Value = -20.06446253974092
$ws.Range("C12").Value = -20.06446253974092
$ws.Range("D12").Value = -20.06446253974092
$ws.Range("E12").Value = -20.06446253974092
$ws.Range("F12").Value = -20.06446253974092
$ws.Range("G12").Value = -20.06446253974092
$ws.Range("H12").Value = -20.06446253974092
$ws.Range("I12").Value = -20.06446253974092
$ws.Range("J12").Value = -20.06446253974092
$ws.Range("K12").Value = -20.06446253974092

$ws.Range("B13").Value = -20.06446253974092
$ws.Range("C13").Value = -20.06446253974092
$ws.Range("D13").Value = -20.06446253974092
$ws.Range("E13").Value = 1.597759461458581
$ws.Range("F13").Value = -20.06446253974092
$ws.Range("G13").Value = -20.06446253974092
$ws.Range("H13").Value = -20.06446253974092
$ws.Range("I13").Value = -20.06446253974092
$ws.Range("J13").Value = 2.269659318285476
$ws.Range("K13").Value = 1.626458329300656

$ws.Range("B14").Value = -20.06446253974092
$ws.Range("C14").Value = -20.06446253974092
$ws.Range("D14").Value = 1.68621051190162
$ws.Range("E14").Value = -20.06446253974092
$ws.Range("F14").Value = -20.06446253974092
$ws.Range("G14").Value = -20.06446253974092
$ws.Range("H14").Value = -20.06446253974092
$ws.Range("I14").Value = -20.06446253974092
$ws.Range("J14").Value = -20.06446253974092
$ws.Range("K14").Value = 2.105905421000831

$ws.Range("B15").Value = -20.06446253974092
$ws.Range("C15").Value = -20.06446253974092
$ws.Range("D15").Value = -0.1779368623017062
$ws.Range("E15").Value = -20.06446253974092
$ws.Range("F15").Value = -20.06446253974092
$ws.Range("G15").Value = -20.06446253974092
$ws.Range("H15").Value = -20.06446253974092
$ws.Range("I15").Value = -20.06446253974092
$ws.Range("J15").Value = -20.06446253974092
$ws.Range("K15").Value = -20.06446253974092

$ws.Range("B16").Value = -20.06446253974092
$ws.Range("C16").Value = -20.06446253974092
$ws.Range("D16").Value = -20.06446253974092
$ws.Range("E16").Value = -20.06446253974092
$ws.Range("F16").Value = -20.06446253974092
$ws.Range("G16").Value = -20.06446253974092
$ws.Range("H16").Value = -20.06446253974092
$ws.Range("I16").Value = -20.06446253974092
$ws.Range("J16").Value = 2.315816986948407
$ws.Range("K16").Value = -20.06446253974092

$ws.Range("B17").Value = -20.06446253974092
$ws.Range("C17").Value = 0.7022709738173515
$ws.Range("D17").Value = -0.3322470754012632
$ws.Range("E17").Value = -20.06446253974092
$ws.Range("F17").Value = -20.06446253974092
$ws.Range("G17").Value = -20.06446253974092
$ws.Range("H17").Value = 0.556799696493311
$ws.Range("I17").Value = 0.9056974295707735
$ws.Range("J17").Value = 1.269566817045639
$ws.Range("K17").Value = -20.06446253974092

$ws.Range("B18").Value = -20.06446253974092
$ws.Range("C18").Value = -20.06446253974092
$ws.Range("D18").Value = -20.06446253974092
$ws.Range("E18").Value = -20.06446253974092
$ws.Range("F18").Value = -20.06446253974092
$ws.Range("G18").Value = -20.06446253974092
$ws.Range("H18").Value = 0.4430195193132901
$ws.Range("I18").Value = 0.9590975886223637
$ws.Range("J18").Value = 1.379009695973175
$ws.Range("K18").Value = -20.06446253974092

$ws.Range("B19").Value = -20.06446253974092
$ws.Range("C19").Value = -20.06446253974092
$ws.Range("D19").Value = 1.577164574682874
$ws.Range("E19").Value = -20.06446253974092
$ws.Range("F19").Value = -20.06446253974092
$ws.Range("G19").Value = -20.06446253974092
$ws.Range("H19").Value = 1.884949923072709
$ws.Range("I19").Value = 2.106565480622832
$ws.Range("J19").Value = -20.06446253974092
$ws.Range("K19").Value = -20.06446253974092

$ws.Range("B20").Value = -20.06446253974092
$ws.Range("C20").Value = 1.675615058126258
$ws.Range("D20").Value = 2.090337340241923
$ws.Range("E20").Value = -20.06446253974092
$ws.Range("F20").Value = 3.850464599864329
$ws.Range("G20").Value = -20.06446253974092
$ws.Range("H20").Value = 2.197234425206382
$ws.Range("I20").Value = 1.909064223371347
$ws.Range("J20").Value = -20.06446253974092
$ws.Range("K20").Value = 2.431845375479662

$ws.Range("B21").Value = -20.06446253974092
$ws.Range("C21").Value = 1.772903995217078
$ws.Range("D21").Value = -20.06446253974092
$ws.Range("E21").Value = 2.44644091631957
$ws.Range("F21").Value = -20.06446253974092
$ws.Range("G21").Value = 3.265459273275737
$ws.Range("H21").Value = 2.415138396420013
$ws.Range("I21").Value = -20.06446253974092
$ws.Range("J21").Value = -20.06446253974092
$ws.Range("K21").Value = -20.06446253974092
